$wb = $excel.ActiveWorkbook

# --- Sheet "展览" ---
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F4").Value = 14
$ws1.Range("F6").Value = 29
$ws1.Range("F7").Value = 861
$ws1.Range("F8").Value = 42
$ws1.Range("F9").Value = 6824
$ws1.Range("G9").Value = 80
$ws1.Range("F10").Value = 46
$ws1.Range("F11").Value = 102
$ws1.Range("F12").Value = 140
$ws1.Range("F13").Value = 6435
$ws1.Range("F14").Value = 126
$ws1.Range("F15").Value = 270
$ws1.Range("F16").Value = 4350
$ws1.Range("F20").Value = 4322
$ws1.Range("F21").Value = 227
$ws1.Range("F22").Value = 231
$ws1.Range("F23").Value = 317
$ws1.Range("F27").Value = 166
$ws1.Range("F31").Value = 69
$ws1.Range("F32").Value = 7881
$ws1.Range("F34").Value = 1339
$ws1.Range("F35").Value = 653
$ws1.Range("F36").Value = 15
$ws1.Range("F39").Value = 1571
$ws1.Range("F41").Value = 909
$ws1.Range("F43").Value = 3942
$ws1.Range("F46").Value = 107
$ws1.Range("F47").Value = 36
$ws1.Range("F49").Value = 1082

# --- Sheet "全部类型" ---
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F7").Value = 14
$ws4.Range("F9").Value = 29
$ws4.Range("F10").Value = 861
$ws4.Range("F11").Value = 42
$ws4.Range("F12").Value = 6824
$ws4.Range("G12").Value = 80
$ws4.Range("F13").Value = 46
$ws4.Range("F14").Value = 102
$ws4.Range("F15").Value = 140
$ws4.Range("F16").Value = 6435
$ws4.Range("F17").Value = 126
$ws4.Range("F18").Value = 270
$ws4.Range("F19").Value = 4350
$ws4.Range("F22").Value = 4322
$ws4.Range("F23").Value = 227
$ws4.Range("F24").Value = 231
$ws4.Range("F25").Value = 317
$ws4.Range("F31").Value = 69
$ws4.Range("F33").Value = 7881
$ws4.Range("F35").Value = 1339
$ws4.Range("F36").Value = 653
$ws4.Range("F39").Value = 1571
$ws4.Range("F41").Value = 909
$ws4.Range("F43").Value = 3942
$ws4.Range("F46").Value = 107
$ws4.Range("F48").Value = 1082
